$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Undo the earlier test edit: remove the BIOL222 row (row 9) that was added,
# shifting all subsequent rows back up.
$ws.Rows.Item(9).Delete()

$ws.Range("E10").Select()
